$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cur = $cell.Value2
    if ($cur -like "*_old") {
        $cell.Value2 = $cur.Substring(0, $cur.Length - 4) + "_FV2210"
    } elseif ($cur -like "*_new") {
        $cell.Value2 = $cur.Substring(0, $cur.Length - 4) + "_FV2304"
    }
}

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U89"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
